$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5482.75
$ws.Range("I7").Value = 2400
$ws.Range("K7").Value = 2400
$ws.Range("M7").Value = -2288
$ws.Range("H14").Value = 5482.75
$ws.Range("I14").Value = 2400
$ws.Range("K14").Value = 2400
$ws.Range("M14").Value = -2209
$ws.Range("H17").Value = 2137.2222
$ws.Range("J17").Value = 2300.2173
$ws.Range("L17").Value = 6900.651899999999
$ws.Range("N17").Value = -7236.651899999999
$ws.Range("H94").Value = 8401.200000000001
$ws.Range("I94").Value = 8401.200000000001
$ws.Range("K94").Value = 8401.200000000001
$ws.Range("M94").Value = -7950.200000000001
$ws.Range("H138").Value = 1824.4706
$ws.Range("I138").Value = 1801.1875
$ws.Range("K138").Value = 5403.5625
$ws.Range("M138").Value = -263.5625

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11614.35
$ws.Range("J44").Value = 11614.35
$ws.Range("L44").Value = 11614.35
$ws.Range("N44").Value = -12590.35
$ws.Range("H55").Value = 58999.5
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 58999.5
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 58999.5
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -59629.5
$ws.Range("H74").Value = 7301.8125
$ws.Range("I74").Value = 6951.923
$ws.Range("K74").Value = 6951.923
$ws.Range("M74").Value = -6077.923
$ws.Range("H77").Value = 7301.8125
$ws.Range("I77").Value = 6951.923
$ws.Range("K77").Value = 34759.615
$ws.Range("M77").Value = -30391.615
$ws.Range("H97").Value = 2179.5715
$ws.Range("I97").Value = 1765
$ws.Range("J97").Value = 2490.5
$ws.Range("K97").Value = 1765
$ws.Range("L97").Value = 2490.5
$ws.Range("M97").Value = -1269
$ws.Range("N97").Value = -3482.5
$ws.Range("H132").Value = 2661.1875
$ws.Range("I132").Value = 1548.3334
$ws.Range("K132").Value = 4645.0002
$ws.Range("M132").Value = -2115.0002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 762.5
$ws.Range("I29").Value = 716.6667
$ws.Range("K29").Value = 716.6667
$ws.Range("M29").Value = -427.6667
$ws.Range("H57").Value = 94996
$ws.Range("J57").Value = 94996
$ws.Range("L57").Value = 94996
$ws.Range("N57").Value = -96436
$ws.Range("H134").Value = 1783.2
$ws.Range("I134").Value = 1783.2
$ws.Range("K134").Value = 5349.6
$ws.Range("M134").Value = -2814.6
$ws.Range("H136").Value = 94996
$ws.Range("J136").Value = 94996
$ws.Range("L136").Value = 94996
$ws.Range("N136").Value = -105196

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5145.1387
$ws.Range("I31").Value = 1980.5
$ws.Range("J31").Value = 8309.777
$ws.Range("K31").Value = 1980.5
$ws.Range("L31").Value = 8309.777
$ws.Range("M31").Value = -1685.5
$ws.Range("N31").Value = -8899.777
$ws.Range("H34").Value = 5145.1387
$ws.Range("I34").Value = 1980.5
$ws.Range("J34").Value = 8309.777
$ws.Range("K34").Value = 1980.5
$ws.Range("L34").Value = 8309.777
$ws.Range("M34").Value = -1778.5
$ws.Range("N34").Value = -8713.777
$ws.Range("H58").Value = 2379.25
$ws.Range("I58").Value = 852.61536
$ws.Range("K58").Value = 852.61536
$ws.Range("M58").Value = -649.61536
$ws.Range("H95").Value = 4688.1665
$ws.Range("J95").Value = 4688.1665
$ws.Range("L95").Value = 4688.1665
$ws.Range("N95").Value = -10180.1665
$ws.Range("H96").Value = 11962.3
$ws.Range("J96").Value = 11962.3
$ws.Range("L96").Value = 11962.3
$ws.Range("N96").Value = -17454.3
$ws.Range("H105").Value = 1343.7142
$ws.Range("I105").Value = 1653
$ws.Range("K105").Value = 1653
$ws.Range("M105").Value = 94
$ws.Range("H136").Value = 2379.25
$ws.Range("I136").Value = 852.61536
$ws.Range("K136").Value = 2557.84608
$ws.Range("M136").Value = -7.846080000000256

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1972.25
$ws.Range("I81").Value = 1463
$ws.Range("K81").Value = 4389
$ws.Range("M81").Value = -3266
$ws.Range("H84").Value = 1972.25
$ws.Range("I84").Value = 1463
$ws.Range("K84").Value = 13167
$ws.Range("M84").Value = -7551
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4827.5454
$ws.Range("I55").Value = 5671.5
$ws.Range("K55").Value = 5671.5
$ws.Range("M55").Value = -5344.5
$ws.Range("H102").Value = 2583.077
$ws.Range("I102").Value = 2583.077
$ws.Range("K102").Value = 2583.077
$ws.Range("M102").Value = -961.0770000000002
$ws.Range("H132").Value = 36631.188
$ws.Range("I132").Value = 48098.523
$ws.Range("J132").Value = 7325.778
$ws.Range("K132").Value = 144295.569
$ws.Range("L132").Value = 21977.334
$ws.Range("M132").Value = -141765.569
$ws.Range("N132").Value = -27037.334

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 1500
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H130").Value = 24248.75
$ws.Range("J130").Value = 24248.75
$ws.Range("L130").Value = 24248.75
$ws.Range("N130").Value = -34288.75
$ws.Range("H132").Value = 3364.7646
$ws.Range("I132").Value = 3450.2144
$ws.Range("K132").Value = 10350.6432
$ws.Range("M132").Value = -7820.643199999999
$ws.Range("H136").Value = 2938.3333
$ws.Range("I136").Value = 2334.2
$ws.Range("J136").Value = 5959
$ws.Range("K136").Value = 7002.599999999999
$ws.Range("L136").Value = 17877
$ws.Range("M136").Value = -4452.599999999999
$ws.Range("N136").Value = -22977

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2180.6316
$ws.Range("I136").Value = 1237.6154
$ws.Range("J136").Value = 4223.8335
$ws.Range("K136").Value = 3712.8462
$ws.Range("L136").Value = 12671.5005
$ws.Range("M136").Value = -1162.8462
$ws.Range("N136").Value = -17771.5005
